# Generate Report for Handback
# This applies the "handback generated" update for the 1cd3b2b1-a014-49e3-8774-1f039a9fdc8e
# source file row (row 6) on the zh-cn and de-de worksheets, plus widening the
# "Error Detail" column (P) on both sheets.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6abdb06af2818ea6ec7bb7b01373c08a06a11a47/e2e/1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/292c8daab3629f71322394886a4dae07ff5fe9c1/e2e/1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.md."

$targetMdDisplay = "1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.md"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6abdb06af2818ea6ec7bb7b01373c08a06a11a47/e2e/1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.md", "", "", $targetMdDisplay)
$wsZh.Range("I6").Font.Underline = 2
$wsZh.Range("I6").Font.Color = 15570276
$wsZh.Range("J6").Value = "1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.ae0b30c5bad459ce728d951f976ff966939c99dc.zh-cn.xlf"
$wsZh.Range("K6").Value = "2016-11-08 23:04:12"
$wsZh.Range("P6").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.166666666666664

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I6"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/6abdb06af2818ea6ec7bb7b01373c08a06a11a47/e2e/1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.md", "", "", $targetMdDisplay)
$wsDe.Range("I6").Font.Underline = 2
$wsDe.Range("I6").Font.Color = 15570276
$wsDe.Range("J6").Value = "1cd3b2b1-a014-49e3-8774-1f039a9fdc8e.ae0b30c5bad459ce728d951f976ff966939c99dc.de-de.xlf"
$wsDe.Range("K6").Value = "2016-11-08 23:04:30"
$wsDe.Range("P6").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.166666666666664
